# Refresh the "cryptos" price/volume table (GitHub Actions scrape update).
# Two rows (19/20 and 50/51) also swapped which coin they describe, so their
# Coin/Link/Price/Volume columns are rewritten in full; every other changed
# row only updates Price (D) and/or Volume(1h) (E).
#
# Numeric-looking Price strings (e.g. "228.50", "1.00", "0.0701") are written
# with a leading apostrophe so Excel keeps them as text (matching the
# worksheet's existing inline-string cells) instead of collapsing them into
# numbers and dropping significant trailing/leading zeros. The apostrophe is
# Excel's standard "force text" quote-prefix marker and is not stored as part
# of the cell value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.695.17'
$ws.Range('E2').Value = '  +0.70%  '

$ws.Range('D3').Value = '1.820.05'
$ws.Range('E3').Value = '  +1.33%  '

$ws.Range('E4').Value = '  +0.10%  '

$ws.Range('D5').Value = '''228.50'
$ws.Range('E5').Value = '  +0.70%  '

$ws.Range('D6').Value = '''0.576'
$ws.Range('E6').Value = '  +3.73%  '

$ws.Range('D7').Value = '''1.00'
$ws.Range('E7').Value = '  +0.18%  '

$ws.Range('D8').Value = '''34.74'
$ws.Range('E8').Value = '  +7.11%  '

$ws.Range('D9').Value = '''0.300'
$ws.Range('E9').Value = '  +1.32%  '

$ws.Range('D10').Value = '''0.0701'
$ws.Range('E10').Value = '  +1.10%  '

$ws.Range('D11').Value = '''0.0954'
$ws.Range('E11').Value = '  +0.49%  '

$ws.Range('D12').Value = '2.088.56'
$ws.Range('E12').Value = '  +1.65%  '

$ws.Range('D13').Value = '''11.46'
$ws.Range('E13').Value = '  +3.24%  '

$ws.Range('D14').Value = '1.841.69'
$ws.Range('E14').Value = '  +1.83%  '

$ws.Range('D15').Value = '''0.645'
$ws.Range('E15').Value = '  +1.85%  '

$ws.Range('D16').Value = '34.716.64'
$ws.Range('E16').Value = '  +0.88%  '

$ws.Range('D17').Value = '''4.33'
$ws.Range('E17').Value = '  +2.13%  '

$ws.Range('D18').Value = '''69.11'
$ws.Range('E18').Value = '  +1.04%  '

$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').Value = '0.0₃0802'
$ws.Range('E19').Value = '  +0.19%  '

$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').Value = '''246.81'
$ws.Range('E20').Value = '  +0.07%  '

$ws.Range('D21').Value = '''11.62'
$ws.Range('E21').Value = '  +4.49%  '

$ws.Range('E22').Value = '  +0.04%  '

$ws.Range('D23').Value = '''4.18'
$ws.Range('E23').Value = '  +0.37%  '

$ws.Range('D24').Value = '''174.12'
$ws.Range('E24').Value = '  +6.94%  '

$ws.Range('E25').Value = '  +1.28%  '

$ws.Range('D26').Value = '''7.49'
$ws.Range('E26').Value = '  +3.25%  '

$ws.Range('D27').Value = '''16.81'
$ws.Range('E27').Value = '  +2.21%  '

$ws.Range('E28').Value = '  +2.22%  '

$ws.Range('E29').Value = '  -0.11%  '

$ws.Range('E30').Value = '  +2.47%  '

$ws.Range('D31').Value = '''0.0532'
$ws.Range('E31').Value = '  +1.65%  '

$ws.Range('D32').Value = '''3.85'
$ws.Range('E32').Value = '  +2.06%  '

$ws.Range('E33').Value = '  +0.87%  '

$ws.Range('D34').Value = '''1.85'
$ws.Range('E34').Value = '  +1.19%  '

$ws.Range('E35').Value = '  +1.06%  '

$ws.Range('D36').Value = '1.412.59'
$ws.Range('E36').Value = '  -2.09%  '

$ws.Range('D37').Value = '''0.683'
$ws.Range('E37').Value = '  +1.97%  '

$ws.Range('D38').Value = '''1.07'
$ws.Range('E38').Value = '  +1.86%  '

$ws.Range('D39').Value = '''0.0192'
$ws.Range('E39').Value = '  +0.44%  '

$ws.Range('D40').Value = '''84.93'
$ws.Range('E40').Value = '  +1.05%  '

$ws.Range('E41').Value = '  +4.56%  '

$ws.Range('D42').Value = '''0.957'
$ws.Range('E42').Value = '  +2.53%  '

$ws.Range('E43').Value = '  -0.10%  '

$ws.Range('D44').Value = '''13.75'
$ws.Range('E44').Value = '  -0.11%  '

$ws.Range('E46').Value = '  -1.16%  '

$ws.Range('D47').Value = '''6.09'
$ws.Range('E47').Value = '  -0.04%  '

$ws.Range('D48').Value = '1.989.14'
$ws.Range('E48').Value = '  +1.98%  '

$ws.Range('D49').Value = '''105.54'
$ws.Range('E49').Value = '  -0.16%  '

$ws.Range('B50').Value = 'PaxDollar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D50').Value = '''1.00'
$ws.Range('E50').Value = '  +0.03%  '

$ws.Range('B51').Value = 'BitcoinSV'
$ws.Range('C51').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D51').Value = '''50.16'
$ws.Range('E51').Value = '  +0.54%  '

